$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.266.64"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.856.74"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.19"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4730"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06429"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").Value = "1.851.12"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07436"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.08"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.988"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.44"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6327"
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("D16").Value = "30.235.54"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.78"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007323"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "225.60"
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("D21").Value = "2.092.32"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.105"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.033"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.38"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.236"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.80"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("E28").Value = "  -5.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1026"
$ws.Range("E29").Value = "  +9.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.379"
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.228"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.902"
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04885"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.148"
$ws.Range("E34").Value = "  -4.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7266"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9994"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01923"
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.624"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9009"
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.976"
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.87"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9937"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4093"
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.531"
$ws.Range("E45").Value = "  -5.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.034"
$ws.Range("E46").Value = "  -5.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.30"
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1204"
$ws.Range("E48").Value = "  -5.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.769"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.402"
$ws.Range("E50").Value = "  -4.74%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.95"
$ws.Range("E51").Value = "  -2.26%  "
